# Docx writer: Use different style for block quotes in notes.
# Adds a new paragraph style "Footnote Block Text" (styleId
# "FootnoteBlockText"), based on / followed-by "Footnote Text",
# mirroring the existing "Block Text" style's paragraph formatting
# (spacing before/after 100, left/right indent 480, no first-line
# indent) so block quotes inside footnotes can get their own look.

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1
$style = $d.Styles.Add("Footnote Block Text", 1)

$style.BaseStyle = "Footnote Text"
$style.NextParagraphStyle = "Footnote Text"
$style.Priority = 9
$style.UnhideWhenUsed = $true
$style.QuickStyle = $true

# ParagraphFormat distances are in points; twips = points * 20, so
# 100 twips -> 5pt, 480 twips -> 24pt.
$style.ParagraphFormat.SpaceBefore = 5
$style.ParagraphFormat.SpaceAfter = 5
$style.ParagraphFormat.FirstLineIndent = 0
$style.ParagraphFormat.LeftIndent = 24
$style.ParagraphFormat.RightIndent = 24
